$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requerimiento Inicial")

# Update requirement names (column B) for RF.1 - RF.14 (rows 2-15)
$ws.Range("B2").Value = "El usuario se registra en el sistema"
$ws.Range("B3").Value = "permitir acceso a usuarios mediante autenticación"
$ws.Range("B4").Value = "El sistema asigna roles a cada usuario"
$ws.Range("B5").Value = "El usuario visualiza el catálogo de productos"
$ws.Range("B6").Value = "El usuario filtra productos por categoría"
$ws.Range("B7").Value = "El usuario agrega productos al carrito de compras"
$ws.Range("B8").Value = "El usuario gestiona las cantidades de productos en `nel carrito"
$ws.Range("B9").Value = "El usuario gestiona sus pedidos realizados"
$ws.Range("B10").Value = "El usuario crea un pedido personalizado"
$ws.Range("B11").Value = "El usuario realiza un pedido con selección de zona `nde despacho"
$ws.Range("B12").Value = "El usuario realiza seguimiento del estado de su `npedido"
$ws.Range("B13").Value = "El administrador gestiona el sistema desde el `nbackoffice"
$ws.Range("B14").Value = "El administrador genera reportes y los exporta a `nPDF y Excel"
$ws.Range("B15").Value = "El sistema envía notificaciones por correo electrónico"

# Column B width (target stored width 50.29; engine quantizes ColumnWidth to
# 1/6-character pixel steps, so 49.5 is the closest input that serializes to
# the nearest achievable stored width)
$ws.Columns.Item(2).ColumnWidth = 49.5

# Row heights
$ws.Rows.Item(7).RowHeight = 33.0
$ws.Rows.Item(8).RowHeight = 28.5
$ws.Rows.Item(11).RowHeight = 29.25
$ws.Rows.Item(12).RowHeight = 28.5
$ws.Rows.Item(13).RowHeight = 27.75
$ws.Rows.Item(14).RowHeight = 28.5
